# Updates cryptocurrency Price (column D) and Volume(1h) (column E)
# values on worksheet row 2-51, matching the scraped-data refresh commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = '62.297.48' }
    @{ Cell = "E2"; Value = '  +0.50%  ' }
    @{ Cell = "D3"; Value = '2.997.10' }
    @{ Cell = "E3"; Value = '  -1.13%  ' }
    @{ Cell = "D4"; Value = '0.999' }
    @{ Cell = "E4"; Value = '  -0.08%  ' }
    @{ Cell = "D5"; Value = '543.58' }
    @{ Cell = "E5"; Value = '  -1.93%  ' }
    @{ Cell = "D6"; Value = '138.40' }
    @{ Cell = "E6"; Value = '  +1.43%  ' }
    @{ Cell = "D7"; Value = '0.999' }
    @{ Cell = "E7"; Value = '  -0.10%  ' }
    @{ Cell = "D8"; Value = '2.995.73' }
    @{ Cell = "E8"; Value = '  -0.98%  ' }
    @{ Cell = "D9"; Value = '0.489' }
    @{ Cell = "D10"; Value = '6.93' }
    @{ Cell = "E10"; Value = '  +13.82%  ' }
    @{ Cell = "E11"; Value = '  -1.06%  ' }
    @{ Cell = "D12"; Value = '0.446' }
    @{ Cell = "E12"; Value = '  -1.51%  ' }
    @{ Cell = "E13"; Value = '  -1.46%  ' }
    @{ Cell = "D14"; Value = '33.92' }
    @{ Cell = "E14"; Value = '  -1.87%  ' }
    @{ Cell = "D15"; Value = '3.470.88' }
    @{ Cell = "E15"; Value = '  -1.34%  ' }
    @{ Cell = "D16"; Value = '62.246.20' }
    @{ Cell = "E16"; Value = '  +0.24%  ' }
    @{ Cell = "D17"; Value = '2.996.95' }
    @{ Cell = "E17"; Value = '  -1.01%  ' }
    @{ Cell = "E18"; Value = '  -2.12%  ' }
    @{ Cell = "D19"; Value = '6.56' }
    @{ Cell = "E19"; Value = '  -1.95%  ' }
    @{ Cell = "D20"; Value = '469.93' }
    @{ Cell = "E20"; Value = '  -1.26%  ' }
    @{ Cell = "D21"; Value = '13.41' }
    @{ Cell = "E21"; Value = '  +0.53%  ' }
    @{ Cell = "D22"; Value = '0.653' }
    @{ Cell = "E22"; Value = '  -3.72%  ' }
    @{ Cell = "D23"; Value = '7.16' }
    @{ Cell = "E23"; Value = '  +0.70%  ' }
    @{ Cell = "D24"; Value = '79.41' }
    @{ Cell = "E24"; Value = '  -1.35%  ' }
    @{ Cell = "D25"; Value = '12.60' }
    @{ Cell = "E25"; Value = '  +3.38%  ' }
    @{ Cell = "E26"; Value = '  -0.23%  ' }
    @{ Cell = "E27"; Value = '  -0.92%  ' }
    @{ Cell = "D28"; Value = '7.63' }
    @{ Cell = "E28"; Value = '  -2.71%  ' }
    @{ Cell = "D29"; Value = '2.01' }
    @{ Cell = "E29"; Value = '  +4.31%  ' }
    @{ Cell = "E30"; Value = '  +0.21%  ' }
    @{ Cell = "D31"; Value = '25.36' }
    @{ Cell = "E31"; Value = '  -2.06%  ' }
    @{ Cell = "E32"; Value = '  -2.98%  ' }
    @{ Cell = "D33"; Value = '2.35' }
    @{ Cell = "E33"; Value = '  +0.88%  ' }
    @{ Cell = "E34"; Value = '  +1.21%  ' }
    @{ Cell = "D35"; Value = '54.68' }
    @{ Cell = "E35"; Value = '  -2.12%  ' }
    @{ Cell = "D36"; Value = '5.84' }
    @{ Cell = "E36"; Value = '  -1.80%  ' }
    @{ Cell = "D37"; Value = '451.76' }
    @{ Cell = "E37"; Value = '  -2.07%  ' }
    @{ Cell = "D38"; Value = '0.0811' }
    @{ Cell = "E38"; Value = '  +1.10%  ' }
    @{ Cell = "D39"; Value = '0.0392' }
    @{ Cell = "E39"; Value = '  +1.28%  ' }
    @{ Cell = "D40"; Value = '2.945.72' }
    @{ Cell = "E40"; Value = '  -8.47%  ' }
    @{ Cell = "E41"; Value = '  -4.08%  ' }
    @{ Cell = "D42"; Value = '8.05' }
    @{ Cell = "E42"; Value = '  -1.55%  ' }
    @{ Cell = "D43"; Value = '2.54' }
    @{ Cell = "E43"; Value = '  +2.38%  ' }
    @{ Cell = "D44"; Value = '26.72' }
    @{ Cell = "E44"; Value = '  +2.62%  ' }
    @{ Cell = "E45"; Value = '  -0.05%  ' }
    @{ Cell = "E46"; Value = '  +1.20%  ' }
    @{ Cell = "D47"; Value = '2.00' }
    @{ Cell = "E47"; Value = '  -0.60%  ' }
    @{ Cell = "E48"; Value = '  +0.36%  ' }
    @{ Cell = "D49"; Value = '114.96' }
    @{ Cell = "E49"; Value = '  -3.08%  ' }
    @{ Cell = "D50"; Value = '0.0₃0493' }
    @{ Cell = "E50"; Value = '  -1.11%  ' }
    @{ Cell = "E51"; Value = '  -1.20%  ' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    # Force text storage (source cells are text, e.g. "1.00", "0.999")
    # instead of letting Excel auto-convert to a number.
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    # Drop the explicit Text number-format again so the cell keeps
    # the workbook default (General) style, same as before the edit.
    $cell.ClearFormats()
}
